$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, `
                                $true, 1, $false, $newText, 2)
    if (-not $found) {
        throw "Could not find text: $oldText"
    }
    Write-Host "Replaced: $($oldText.Substring(0, [Math]::Min(40, $oldText.Length)))..."
}

# 1. Professional summary paragraph
$oldSummary = "Experienced back-end developer and DevOps engineer with 15+ years in IT background and 11 years in development. My account is in top 1% of all Ruby accounts on Github (according to GitAwards). I am a programmer polyglot always eager to learn new technologies. Most recently I've been working with Ruby and JS and now curious about GO. Problem solver, result-oriented and self-starter. Comfortable working remotely with teammates."
$newSummary = "As an experienced Site Reliability Engineer and Back-End Engineer with over 20 years in IT development, I bring a wealth of knowledge and expertise to any project. I am a certified Google Cloud Architect and am passionate about working with Kubernetes and automating tasks to streamline processes. My Github account is in the top 1% of all Ruby accounts on GitAwards, and I am a polyglot programmer who is always eager to learn new technologies. While my recent focus has been on Ruby and JS, I am also curious about exploring Go. As a problem solver, I am highly result-oriented and self-motivated, and I am comfortable working remotely with teammates. I believe in delivering quality work that exceeds expectations, and I am excited about the opportunity to bring my skills to your organization."
Replace-Text $oldSummary $newSummary

# 2. "Working on GCP cloud..." body text paragraph
$oldGcp = "Working on GCP cloud, with k8s (GKE) on high load site with a few dozens of company and external services."
$newGcp = "I am working on GCP cloud, with k8s (GKE) on a high-load site with a few dozen company and external services."
Replace-Text $oldGcp $newGcp

# 3. "Moved CI/CD to GCP cloudBuild" bullet point
$oldCicd = "Moved CI/CD to GCP cloudBuild"
$newCicd = "Moved (or created from scratch) CI/CD to GCP cloudBuild from different providers"
Replace-Text $oldCicd $newCicd

# 4. Add new bullet "WIP on MySQL cluster seamless upgrade" after the
#    "Optimised some high load ruby code..." bullet item.
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*Optimised some high load ruby code to achieve 100x more CPU efficiency*") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not find target bullet paragraph"
}
$targetPara = $d.Paragraphs($targetIndex)
$targetPara.Range.InsertParagraphAfter()
$newBulletRange = $d.Paragraphs($targetIndex + 1).Range
$newBulletRange.MoveEnd(1, -1)
$newBulletRange.Text = "WIP on MySQL cluster seamless upgrade"
Write-Host "Inserted new bullet paragraph at index $($targetIndex + 1)"
